$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 8.742999999999999
$ws.Range("D4").Value = -7.968999999999999
$ws.Range("B6").Value = 6.279000000000001
$ws.Range("B7").Value = 5.668
$ws.Range("B8").Value = 6.01
$ws.Range("D8").Value = -7.965999999999999
$ws.Range("D9").Value = -7.578
$ws.Range("D12").Value = -7.072
$ws.Range("B16").Value = 4.947000000000001
$ws.Range("D17").Value = -8.268999999999998
$ws.Range("D18").Value = -8.463000000000001
$ws.Range("D19").Value = -7.967000000000001
$ws.Range("B20").Value = 8.24
$ws.Range("D20").Value = -7.853
$ws.Range("B21").Value = 9.094000000000001
$ws.Range("D26").Value = -7.641999999999999
$ws.Range("B28").Value = 5.497999999999999
$ws.Range("B29").Value = 5.281
$ws.Range("B30").Value = 5.209000000000001
$ws.Range("D31").Value = -7.914
$ws.Range("B32").Value = 6.48
$ws.Range("D39").Value = -7.632
$ws.Range("B40").Value = 9.435999999999998
$ws.Range("D40").Value = -8.019
$ws.Range("D41").Value = -7.998
$ws.Range("D42").Value = -8.036
$ws.Range("D43").Value = -7.777000000000001
$ws.Range("B46").Value = 5.613
$ws.Range("D47").Value = -7.639
$ws.Range("D48").Value = -7.651999999999999
$ws.Range("B51").Value = 5.377999999999999
$ws.Range("B52").Value = 5.776
$ws.Range("D54").Value = -8.135
$ws.Range("B57").Value = 5.040999999999999
$ws.Range("B59").Value = 5.298
$ws.Range("B62").Value = 5.395999999999999
$ws.Range("D62").Value = -8.252000000000001
$ws.Range("D63").Value = -7.253
$ws.Range("D64").Value = -7.263
$ws.Range("B66").Value = 5.211
$ws.Range("B73").Value = 7.224000000000001
$ws.Range("B74").Value = 9.165000000000001
$ws.Range("D76").Value = -7.748
$ws.Range("B77").Value = 6.241000000000001
$ws.Range("D81").Value = -8.109999999999999
$ws.Range("D84").Value = -8.372
$ws.Range("D89").Value = -8.260999999999999
$ws.Range("B92").Value = 5.178
$ws.Range("D94").Value = -7.640000000000001
$ws.Range("B100").Value = 6.383
